$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "intersection"/"columnName" lookup-table rows appended under the
# existing BC1:BD4 ("intersection"/"columnName") table.
# (Column order chosen to reproduce the exact shared-string insertion
# order recorded in the saved workbook.)
$ws.Range("BC5").Value = "Socio economic class"
$ws.Range("BD5").Value = "Q10"

$ws.Range("BC6").Value = "Religion"
$ws.Range("BD6").Value = "Q9"

$ws.Range("BC7").Value = "Dietary restrictions"

$ws.Range("BC8").Value = "Class"
$ws.Range("BD8").Value = "Q12"

# This entry was (per author's commit) placed two columns further right
# than the rest of the table (BF/BG instead of BC/BD).
$ws.Range("BF13").Value = "Lunch period"
$ws.Range("BG13").Value = "Q32"

$ws.Range("BD7").Value = "test"

# Header row: Z1 used to reference "Q10" (shared string), now should read "test"
$ws.Range("Z1").Value = "test"

# Restore the view state captured in the saved workbook.
$excel.ActiveWindow.Zoom = 58
$ws.Range("BC8:BD8").Select()
